# Fruta / hortaliza, semanal
# A new weekly observation is inserted as row 191 (shifting the existing
# rows 191-222 down to 192-223), matching the OOXML diff where every row
# from 191 onward is pushed down by one and a brand-new record appears
# at row 191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 191; Excel shifts rows 191..222
# down to 192..223 and copies formatting (incl. the date style) from the
# row above, which already matches the style used throughout column D.
$ws.Rows.Item(191).EntireRow.Insert()

# Populate the newly inserted row 191 with the new record's data.
$ws.Range("A191").Value2 = 7
$ws.Range("B191").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C191").Value2 = "Ñuble"
$ws.Range("D191").Value2 = 44644
$ws.Range("E191").Value2 = 16
$ws.Range("F191").Value2 = 100112009
$ws.Range("G191").Value2 = "Acelga"
$ws.Range("H191").Value2 = "Sin especificar"
$ws.Range("I191").Value2 = "Primera"
$ws.Range("J191").Value2 = 60
$ws.Range("K191").Value2 = 550
$ws.Range("L191").Value2 = 600
$ws.Range("M191").Value2 = 575
$ws.Range("N191").Value2 = "$/atado 0,5 a 1 kilo"
$ws.Range("O191").Value2 = "Provincia de Diguillín"
$ws.Range("P191").Value2 = 575
$ws.Range("Q191").Value2 = 1
$ws.Range("R191").Value2 = "Hortaliza"
